$p = $ppt.ActivePresentation

# --- Slide 1: remove the "Subtitle 2" placeholder shape (date/time subtitle
#     line under the title) entirely. The deck's title-slide layout defines a
#     subTitle placeholder, so after removing the shape once the host may
#     re-synthesize an empty placeholder shape; keep deleting any remaining
#     subtitle placeholder until the slide truly has none left. ---
$slide1 = $p.Slides.Item(1)

$maxTries = 5
for ($try = 1; $try -le $maxTries; $try++) {
    $subtitleShape = $null
    for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
        $shp = $slide1.Shapes.Item($i)
        if ($shp.PlaceholderFormat.Type -eq 4) {
            $subtitleShape = $shp
        }
    }
    if ($subtitleShape -eq $null) {
        break
    }
    $subtitleShape.Delete()
}

# --- Slide 8: the "edit/run src/Plasmodium.sh" bullet used to be split across
#     two separate runs ("src/" and "Plasmodium.sh") sharing identical run
#     properties; merge them into a single run "src/Plasmodium.sh" while
#     leaving every other run/paragraph in the shape untouched. ---
$slide8 = $p.Slides.Item(8)

$contentShape = $null
for ($i = 1; $i -le $slide8.Shapes.Count; $i++) {
    $shp = $slide8.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text.Contains("Plasmodium.sh")) {
        $contentShape = $shp
    }
}

$fullRange = $contentShape.TextFrame.TextRange
$fullText = $fullRange.Text
$needle = "src/Plasmodium.sh"
$startPos = $fullText.IndexOf($needle)
if ($startPos -ge 0) {
    $mergedRun = $fullRange.Characters($startPos + 1, $needle.Length)
    $mergedRun.Text = $needle
}
